$d = $word.ActiveDocument

# Locate the last paragraph of the document body ("D:\OC\P8>") and
# create a fresh empty paragraph right after it to serve as the anchor
# for the new content.
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($count)
$anchor = $lastPara.Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()

$count2 = $d.Paragraphs.Count
$targetPara = $d.Paragraphs($count2)
$targetRange = $targetPara.Range

$innerXml = '<w:p><w:proofErr w:type="gramStart"/><w:r><w:t>git</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>branch</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> -M main</w:t></w:r><w:r><w:t xml:space="preserve"> == Pour créer la branche principale </w:t></w:r></w:p><w:p><w:proofErr w:type="gramStart"/><w:r><w:t>git</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> push -u </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>origin</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>main</w:t></w:r><w:r><w:t xml:space="preserve">  =</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">= je pousse sur cette branche </w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Ajouter un nouveau fichier </w:t></w:r></w:p><w:p><w:proofErr w:type="gramStart"/><w:r><w:t>cd</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> D:\OC\P8</w:t></w:r></w:p><w:p><w:proofErr w:type="gramStart"/><w:r><w:t>git</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>add</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Comprendre_08.docx</w:t></w:r></w:p><w:p><w:proofErr w:type="gramStart"/><w:r><w:t>git</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> commit -m "Ajout du fichier Word Comprendre_08"</w:t></w:r></w:p><w:p><w:proofErr w:type="gramStart"/><w:r><w:t>git</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> push</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Mettre </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>a</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> jour les modifications </w:t></w:r></w:p><w:p><w:proofErr w:type="gramStart"/><w:r><w:t>git</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>add</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Comprendre_08.docx</w:t></w:r></w:p><w:p><w:proofErr w:type="gramStart"/><w:r><w:t>git</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> commit -m "Mise à jour du fichier Word Comprendre_08"</w:t></w:r></w:p><w:p><w:proofErr w:type="gramStart"/><w:r><w:t>git</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> push</w:t></w:r></w:p>'

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + '<w:body>' + $innerXml + '</w:body>' + '</w:document></pkg:xmlData></pkg:part></pkg:package>'

$targetRange.InsertXML($packageXml)

# InsertXML leaves one stray empty paragraph mark at the very end of the
# document (the original trailing paragraph mark that the body must always
# keep). Merge it away so the document ends exactly on the last inserted
# paragraph ("git push"), matching the source transcript.
$finalCount = $d.Paragraphs.Count
$secondLast = $d.Paragraphs($finalCount - 1)
$mergePoint = $secondLast.Range.End
$markRange = $d.Range($mergePoint - 1, $mergePoint)
$markRange.Delete()

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
